$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New red fill color (matches the new 4th fill / 3rd cellXf added to styles.xml)
$ws.Range("A2:V2").Interior.ColorIndex = 3

# Special "Wingdings"-style private-use-area bullet character used for the
# third level of the Table of Content outline (U+F0A7).
$sq = [char]0xF0A7

$ws.Range("A2").Value = "Technical Ceramics Market Size, Share, Growth Analysis, By Raw Material (Alumina Ceramics, Titanate Ceramics), By Product (Monolithic ceramics, Ceramic coatings), By Application, By End-use, By Region - Industry Forecast 2025-2032"
$ws.Range("B2").Value = "SQMIG15E2424"
$ws.Range("C2").Value = "https://www.skyquestt.com/report/technical-ceramics-market"
$ws.Range("E2").Value = "157"

$ws.Range("J2").Value = "Report details not available."

$ws.Range("K2").Value = @"
• Introduction
o Objectives of the Study
o Scope of the Report
o Definitions
• Research Methodology
o Information Procurement
o Secondary & Primary Data Methods
o Market Size Estimation
o Market Assumptions & Limitations
• Executive Summary
o Global Market Outlook
o Supply & Demand Trend Analysis
o Segmental Opportunity Analysis
• Market Dynamics & Outlook
o Market Overview
o Market Size
o Market Dynamics
$sq Drivers & Opportunities
$sq Restraints & Challenges
o Porters Analysis
$sq Competitive rivalry
$sq Threat of substitute
$sq Bargaining power of buyers
$sq Threat of new entrants
$sq Bargaining power of suppliers
• Key Market Insights
o Key Success Factors
o Degree of Competition
o Top Investment Pockets
o Market Ecosystem
o Market Attractiveness Index, 2024
o PESTEL Analysis
o Macro-Economic Indicators
o Value Chain Analysis
o Pricing Analysis
o Raw Material Analysis
o Technology Analysis
o Patent Analysis
o Raw Material Analysis
• Global Technical Ceramics Market Size by Raw Material & CAGR (2025-2032)
o Market Overview
o Alumina Ceramics
o Titanate Ceramics
o Zirconate Ceramics
o Ferrite Ceramics
o Aluminium Nitride
o Silicon Carbide
o Silicon Nitride
o Others
• Global Technical Ceramics Market Size by Product & CAGR (2025-2032)
o Market Overview
o Monolithic ceramics
o Ceramic coatings
o Ceramic matrix composites
o Others
• Global Technical Ceramics Market Size by Application & CAGR (2025-2032)
o Market Overview
o Electrical equipment
o Catalyst supports
o Electronic devices
o Wear parts
o Engine parts
o Filters
o Bioceramics
o Others
• Global Technical Ceramics Market Size by End-use & CAGR (2025-2032)
o Market Overview
o Electrical & Electronics
o Automotive
o Machinery
o Environmental
o Medical
o Military & Defense
o Others
• Global Technical Ceramics Market Size & CAGR (2025-2032)
o North America (Raw Material, Product, Application, End-use)
$sq US
$sq Canada
o Europe (Raw Material, Product, Application, End-use)
$sq Germany
$sq Spain
$sq France
$sq UK
$sq Italy
$sq Rest of Europe
o Asia Pacific (Raw Material, Product, Application, End-use)
$sq China
$sq India
$sq Japan
$sq South Korea
$sq Rest of Asia-Pacific
o Latin America (Raw Material, Product, Application, End-use)
$sq Brazil
$sq Rest of Latin America
o Middle East & Africa (Raw Material, Product, Application, End-use)
$sq GCC Countries
$sq South Africa
$sq Rest of Middle East & Africa
• Competitive Intelligence
o Top 5 Player Comparison
o Market Positioning of Key Players, 2024
o Strategies Adopted by Key Market Players
o Recent Developments in the Market
o Company Market Share Analysis, 2024
o Company Profiles of All Key Players
$sq Company Details
$sq Product Portfolio Analysis
$sq Company's Segmental Share Analysis
$sq Revenue Y-O-Y Comparison (2022-2024)
• Key Company Profiles
o Kyocera Corporation (Japan)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
o Morgan Advanced Materials plc (United Kingdom)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
o Schott AG (Germany)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
o CeramTec GmbH (Germany)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
o CoorsTek, Inc. (United States)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
o NGK Insulators Ltd. (Japan)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
o 3M Company (United States)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
o Saint-Gobain S.A. (France)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
o Rauschert GmbH (Germany)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
o Murata Manufacturing Co., Ltd. (Japan)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
o Corning Incorporated (United States)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
o McDanel Advanced Ceramic Technologies LLC (United States)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
o Ceradyne, Inc. (United States)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
o Superior Technical Ceramics Corp. (United States)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
o Krosaki Harima Corporation (Japan)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
o IBIDEN Co., Ltd. (Japan)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
o L3Harris Technologies, Inc. (United States)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
o RHI Magnesita N.V. (Austria)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
o Vesuvius plc (United Kingdom)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
o Sumitomo Electric Industries, Ltd. (Japan)
$sq Company Overview
$sq Business Segment Overview
$sq Financial Updates
$sq Key Developments
• Conclusion & Recommendations
"@

$ws.Range("N2").Value = "Materials"

$ws.Range("P2").Value = @"
◦ Kyocera Corporation (Japan)
◦ Morgan Advanced Materials plc (United Kingdom)
◦ Schott AG (Germany)
◦ CeramTec GmbH (Germany)
◦ CoorsTek, Inc. (United States)
◦ NGK Insulators Ltd. (Japan)
◦ 3M Company (United States)
◦ Saint-Gobain S.A. (France)
◦ Rauschert GmbH (Germany)
◦ Murata Manufacturing Co., Ltd. (Japan)
◦ Corning Incorporated (United States)
◦ McDanel Advanced Ceramic Technologies LLC (United States)
◦ Ceradyne, Inc. (United States)
◦ Superior Technical Ceramics Corp. (United States)
◦ Krosaki Harima Corporation (Japan)
◦ IBIDEN Co., Ltd. (Japan)
◦ L3Harris Technologies, Inc. (United States)
◦ RHI Magnesita N.V. (Austria)
◦ Vesuvius plc (United Kingdom)
◦ Sumitomo Electric Industries, Ltd. (Japan)
"@

$ws.Range("Q2").Value = @"
By Raw Material (Alumina Ceramics, Titanate Ceramics, Zirconate Ceramics, Ferrite Ceramics, Aluminium Nitride, Silicon Carbide, Silicon Nitride, Others), By Product (Monolithic ceramics, Ceramic coatings, Ceramic matrix composites, Others), By Application (Electrical equipment, Catalyst supports, Electronic devices, Wear parts, Engine parts, Filters, Bioceramics, Others), By End-use (Electrical & Electronics, Automotive, Machinery, Environmental, Medical, Military & Defense, Others)
"@

$ws.Range("R2").Value = "5.89"
$ws.Range("S2").Value = "6.27"
$ws.Range("T2").Value = "10.29"
$ws.Range("U2").Value = "6.4%"
